# Automatic update of files.
# Applies the recorded field edits to rows 75-81 of the "Artfynd" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 75 ---
$ws.Range("B75").Value = 96720

# --- Row 76 ---
$ws.Range("A76").Value = 111949575
$ws.Range("B76").Value = 96720
$ws.Range("D76").Value = "VU"
$ws.Range("E76").Value = 220787
$ws.Range("F76").Value = "Knärot"
$ws.Range("G76").Value = "Goodyera repens"
$ws.Range("H76").Value = "(L.) R. Br."
$ws.Range("I76").Value = "'15"
$ws.Range("J76").Value = "plantor/tuvor"
$ws.Range("Q76").Value = 580471
$ws.Range("R76").Value = 7053333
$ws.Range("S76").Value = 1
$ws.Range("Z76").Value = "19:05"
$ws.Range("AB76").Value = "19:05"
$ws.Range("AW76").Value = "Kamilla Andersson"
$ws.Range("AX76").Value = "Kamilla Andersson"

# --- Row 77 ---
$ws.Range("B77").Value = 56575

# --- Row 78 ---
$ws.Range("A78").Value = 111949317
$ws.Range("B78").Value = 96637
$ws.Range("D78").Value = "LC"
$ws.Range("E78").Value = 219790
$ws.Range("F78").Value = "Fläcknycklar"
$ws.Range("G78").Value = "Dactylorhiza maculata"
$ws.Range("H78").Value = "(L.) Soó"
$ws.Range("I78").Value = ""
$ws.Range("Q78").Value = 580500
$ws.Range("R78").Value = 7053329
$ws.Range("Z78").Value = "18:54"
$ws.Range("AB78").Value = "18:54"

# --- Row 79 ---
$ws.Range("A79").Value = 111949678
$ws.Range("B79").Value = 96720
$ws.Range("I79").Value = "'7"
$ws.Range("J79").Value = ""
$ws.Range("Q79").Value = 580467
$ws.Range("R79").Value = 7053330
$ws.Range("S79").Value = 2
$ws.Range("Z79").Value = "19:11"
$ws.Range("AB79").Value = "19:11"
$ws.Range("AW79").Value = "Kim Hultgren"
$ws.Range("AX79").Value = "Kim Hultgren"

# --- Row 80 ---
$ws.Range("B80").Value = 56430

# --- Row 81 ---
$ws.Range("B81").Value = 90794
